$d = $word.ActiveDocument

$replacements = @(
    @("908÷9=", "521÷7="),
    @("618÷9=", "498÷8="),
    @("443÷4=", "268÷4="),
    @("607÷5=", "597÷3="),
    @("995÷5=", "923÷3="),
    @("542÷9=", "319÷5="),
    @("324÷5=", "119÷5="),
    @("866÷5=", "312÷7="),
    @("167÷7=", "624÷8="),
    @("400÷8=", "742÷6="),
    @("649÷7=", "273÷9="),
    @("301÷8=", "923÷4="),
    @("185÷8=", "554÷9="),
    @("430÷6=", "361÷2="),
    @("553÷6=", "457÷9="),
    @("306÷6=", "376÷8="),
    @("584÷4=", "774÷9="),
    @("555÷4=", "249÷4="),
    @("733÷4=", "138÷9="),
    @("996÷4=", "262÷2="),
    @("253÷2=", "557÷6="),
    @("371÷6=", "844÷6="),
    @("938÷2=", "811÷9="),
    @("122÷2=", "777÷7="),
    @("423÷2=", "801÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
